$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B1 from "Response" to "Answer"
$ws.Range("B1").Value = "Answer"

# Select cell B2 (matches the <selection activeCell="B2" sqref="B2"/> left in the saved view)
$ws.Range("B2").Select()
